$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3 values: 24 -> 16 (B3:F3), 14 -> 8 (G3)
$ws.Range("B3:F3").Value = 16
$ws.Range("G3").Value = 8

# Row 4 values: 6 -> 4 (B4:F4), 4 -> 2 (G4)
$ws.Range("B4:F4").Value = 4
$ws.Range("G4").Value = 2

# C8: 24 -> 72
$ws.Range("C8").Value = 72

# Update selection to reflect last user click at F4
$ws.Range("F4").Select() | Out-Null
